$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3170422.8
$ws.Range("I28").Value = 3800107.5
$ws.Range("K28").Value = 3800107.5
$ws.Range("M28").Value = -3799622.5
$ws.Range("H33").Value = 71.46154
$ws.Range("I33").Value = 55.4
$ws.Range("K33").Value = 55.4
$ws.Range("M33").Value = 173.6
$ws.Range("H80").Value = 903.7368
$ws.Range("I80").Value = 1272.3
$ws.Range("J80").Value = 494.22223
$ws.Range("K80").Value = 3816.9
$ws.Range("L80").Value = 1482.66669
$ws.Range("M80").Value = -2818.9
$ws.Range("N80").Value = -3478.66669
$ws.Range("H83").Value = 903.7368
$ws.Range("I83").Value = 1272.3
$ws.Range("J83").Value = 494.22223
$ws.Range("K83").Value = 11450.7
$ws.Range("L83").Value = 4448.00007
$ws.Range("M83").Value = -6458.699999999999
$ws.Range("N83").Value = -14432.00007
$ws.Range("H96").Value = 1962
$ws.Range("I96").Value = 2763
$ws.Range("J96").Value = 360
$ws.Range("K96").Value = 8289
$ws.Range("L96").Value = 1080
$ws.Range("M96").Value = -6916
$ws.Range("N96").Value = -3826
$ws.Range("H98").Value = 2806.6428
$ws.Range("I98").Value = 1774.4166
$ws.Range("J98").Value = 9000
$ws.Range("K98").Value = 1774.4166
$ws.Range("L98").Value = 9000
$ws.Range("M98").Value = -276.4166
$ws.Range("N98").Value = -11996
$ws.Range("H112").Value = 3174.3157
$ws.Range("I112").Value = 399.5
$ws.Range("J112").Value = 3500.7646
$ws.Range("K112").Value = 1198.5
$ws.Range("L112").Value = 10502.2938
$ws.Range("M112").Value = -90.5
$ws.Range("N112").Value = -12718.2938
$ws.Range("H122").Value = 2806.6428
$ws.Range("I122").Value = 1774.4166
$ws.Range("J122").Value = 9000
$ws.Range("K122").Value = 5323.2498
$ws.Range("L122").Value = 27000
$ws.Range("M122").Value = -2873.2498
$ws.Range("N122").Value = -31900
$ws.Range("H132").Value = 1149.3636
$ws.Range("I132").Value = 1154.0312
$ws.Range("K132").Value = 3462.0936
$ws.Range("M132").Value = -932.0935999999997
$ws.Range("H138").Value = 2647.4333
$ws.Range("I138").Value = 2950.9473
$ws.Range("K138").Value = 8852.841899999999
$ws.Range("M138").Value = -3712.841899999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3179.103
$ws.Range("I32").Value = 2048.8545
$ws.Range("K32").Value = 2048.8545
$ws.Range("M32").Value = -1761.8545
$ws.Range("H61").Value = 1968.5555
$ws.Range("I61").Value = 1318.8788
$ws.Range("K61").Value = 1318.8788
$ws.Range("M61").Value = -1106.8788
$ws.Range("H63").Value = 5899.5
$ws.Range("I63").Value = 5899.5
$ws.Range("K63").Value = 5899.5
$ws.Range("M63").Value = -5213.5
$ws.Range("H66").Value = 5899.5
$ws.Range("I66").Value = 5899.5
$ws.Range("K66").Value = 29497.5
$ws.Range("M66").Value = -26065.5
$ws.Range("H97").Value = 455
$ws.Range("I97").Value = 455
$ws.Range("K97").Value = 455
$ws.Range("M97").Value = 41
$ws.Range("H110").Value = 1565.6
$ws.Range("J110").Value = 2030.375
$ws.Range("L110").Value = 2030.375
$ws.Range("N110").Value = -6120.375
$ws.Range("H122").Value = 1814.9445
$ws.Range("I122").Value = 1846.7059
$ws.Range("J122").Value = 1275
$ws.Range("K122").Value = 5540.1177
$ws.Range("L122").Value = 3825
$ws.Range("M122").Value = -3090.1177
$ws.Range("N122").Value = -8725
$ws.Range("H136").Value = 1968.5555
$ws.Range("I136").Value = 1318.8788
$ws.Range("K136").Value = 3956.6364
$ws.Range("M136").Value = -1406.6364

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2888.28
$ws.Range("I31").Value = 1322.7222
$ws.Range("K31").Value = 1322.7222
$ws.Range("M31").Value = -1027.7222
$ws.Range("H34").Value = 2888.28
$ws.Range("I34").Value = 1322.7222
$ws.Range("K34").Value = 1322.7222
$ws.Range("M34").Value = -1120.7222
$ws.Range("H45").Value = 5658.2
$ws.Range("I45").Value = 5658.2
$ws.Range("K45").Value = 5658.2
$ws.Range("M45").Value = -5065.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 514.5714
$ws.Range("J107").Value = 534.1539
$ws.Range("L107").Value = 1602.4617
$ws.Range("N107").Value = -5442.4617
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("N111").ClearContents()
$ws.Range("H131").Value = 5216379
$ws.Range("I131").Value = 83333920
$ws.Range("J131").Value = 8543.111000000001
$ws.Range("K131").Value = 250001760
$ws.Range("L131").Value = 25629.333
$ws.Range("M131").Value = -249996720
$ws.Range("N131").Value = -35709.333

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5749.6665
$ws.Range("H73").Value = 5749.6665
$ws.Range("H122").Value = 2037.9048
$ws.Range("I122").Value = 1840.5333
$ws.Range("K122").Value = 5521.5999
$ws.Range("M122").Value = -3071.5999
$ws.Range("H132").Value = 3539.8572
$ws.Range("I132").Value = 2942.5
$ws.Range("J132").Value = 5265.5557
$ws.Range("K132").Value = 8827.5
$ws.Range("L132").Value = 15796.6671
$ws.Range("M132").Value = -6297.5
$ws.Range("N132").Value = -20856.6671

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12398.111
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 1389.6
$ws.Range("I46").Value = 1299.5
$ws.Range("K46").Value = 1299.5
$ws.Range("M46").Value = -1111.5
$ws.Range("H98").Value = 49899.668
$ws.Range("J98").Value = 49899.668
$ws.Range("L98").Value = 49899.668
$ws.Range("N98").Value = -55889.668

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 60014
$ws.Range("J26").Value = 60014
$ws.Range("L26").Value = 60014
$ws.Range("N26").Value = -60600
$ws.Range("H41").Value = 12833
$ws.Range("J41").Value = 12833
$ws.Range("L41").Value = 12833
$ws.Range("N41").Value = -13613
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H122").Value = 47543.47
$ws.Range("I122").Value = 79285.2
$ws.Range("K122").Value = 237855.6
$ws.Range("M122").Value = -235405.6
$ws.Range("H136").Value = 4010.0435
$ws.Range("I136").Value = 4171.923
$ws.Range("J136").Value = 3799.6
$ws.Range("K136").Value = 12515.769
$ws.Range("L136").Value = 11398.8
$ws.Range("M136").Value = -9965.769
$ws.Range("N136").Value = -16498.8
